$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "246.69"
Set-TextValue $ws.Range("D3") "21.93"
Set-TextValue $ws.Range("D4") "5.409"
Set-TextValue $ws.Range("D5") "0.05790"
Set-TextValue $ws.Range("D6") "3.384"
Set-TextValue $ws.Range("D7") "6.335"
Set-TextValue $ws.Range("D8") "0.8082"
Set-TextValue $ws.Range("D9") "0.9443"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D10") "0.01118"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1428"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.07502"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D13") "0.03197"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D14") "0.03024"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D15") "4.165"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D16") "0.09412"
$ws.Range("E16").Value = "15BitMartTokenBMX"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D17") "0.001586"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D18") "0.04814"
$ws.Range("E18").Value = "17CoinExTokenCET"
Set-TextValue $ws.Range("D19") "0.006182"
Set-TextValue $ws.Range("D20") "0.004113"
Set-TextValue $ws.Range("D21") "0.0009957"
Set-TextValue $ws.Range("D23") "3.772"
Set-TextValue $ws.Range("D25") "0.3231"
Set-TextValue $ws.Range("D27") "0.0003040"
$ws.Range("E27").Value = "26UpBotsUBXT"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.006350"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1074"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.003000"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextValue $ws.Range("D44") "0.006702"
Set-TextValue $ws.Range("D47") "0.3800"
Set-TextValue $ws.Range("D48") "0.1470"
